$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "I/O" mini-table next to the main task table (rows 2-5, cols N/O)
$ws.Range("N2").Value = "I/O"
$ws.Range("O2").Value = "GPIO"
$ws.Range("N3").Value = "Muestreo señal"
$ws.Range("O3").Value = "ADC1_0 (GPIO 36)"
$ws.Range("N4").Value = "NTC Placa caliente"
$ws.Range("O4").Value = "ADC1_3 (GPIO 39)"
$ws.Range("N5").Value = "NTC Placa fría"

# New small table further down the sheet (rows 28-29)
$ws.Range("A28").Value = "I/O"
$ws.Range("B28").Value = "GPIO"
$ws.Range("A29").Value = "Muestreo V"

# Column width adjustments
$ws.Columns.Item(1).ColumnWidth = 14.428571
$ws.Columns.Item(14).ColumnWidth = 16.428571
$ws.Columns.Item(15).ColumnWidth = 16.928571

# Reset the saved scroll position back to the top-left (A1) instead of F1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
